# ARC6-2016.pptx - "Small modification on the slides ARC 6"
#
# Changes applied (per the canonical OOXML diff):
#  Slide 2 ("Grupo 55" group):
#   - ZoneTexte 26 ("Heterogeneous data sources known in advance"): merge the
#     " " and "data sources " runs into a single " data sources " run.
#   - ZoneTexte 28 ("Exported schemata"): drop the paragraph's endParaRPr.
#   - ZoneTexte 29 ("Global schema"): drop the paragraph's endParaRPr.
#  Slide 2 (top level):
#   - ZoneTexte 25: merge "Domenig & Dittrich 1999 Sigmod " and "Record)" runs
#     into a single run.
#  Slide 3 ("Grupo 3" group):
#   - ZoneTexte 26 ("Distributed data services"): drop the endParaRPr.
#   - ZoneTexte 28 ("Exported API"): drop the endParaRPr.
#  Slide 3:
#   - Remove the two click-triggered Fade entrance effects (on "Grouper 25"
#     and "Grouper 30").

$p = $ppt.ActivePresentation

function Get-ShapeByText($container, [string]$text) {
    for ($i = 1; $i -le $container.Count; $i++) {
        $sh = $container.Item($i)
        if ($sh.HasTextFrame -and $sh.TextFrame.HasText) {
            if ($sh.TextFrame.TextRange.Text -eq $text) {
                return $sh
            }
        }
    }
    return $null
}

# Re-types a shape's text so the run that used to carry a (now redundant)
# endParaRPr is rebuilt fresh, without a trailing endParaRPr.
function Remove-EndParaRPr($shape) {
    $text = $shape.TextFrame.TextRange.Text
    $shape.TextFrame.DeleteText()
    $shape.TextFrame.TextRange.Text = $text
}

# ----------------------------------------------------------------------
# Slide 2
# ----------------------------------------------------------------------
$slide2 = $p.Slides.Item(2)
$grp2 = Get-ShapeByText $slide2.Shapes "Data integration: existing work"
$grp2 = $null
for ($i = 1; $i -le $slide2.Shapes.Count; $i++) {
    $sh = $slide2.Shapes.Item($i)
    if ($sh.Type -eq 6 -and $sh.Name -eq "Grupo 55") {
        $grp2 = $sh
    }
}

# -- "Heterogeneous" + " " + "data sources " + "known in advance" --------
$heteroShape = Get-ShapeByText $grp2.GroupItems "Heterogeneous data sources known in advance"
$tr = $heteroShape.TextFrame2.TextRange
# "Heterogeneous" is 13 characters; " data sources " (14 chars) follows it.
$merged = $tr.Characters(14, 14)
$merged.Text = " data sources "

# -- "Exported schemata" --------------------------------------------------
$exportedSchemataShape = Get-ShapeByText $grp2.GroupItems "Exported schemata"
Remove-EndParaRPr $exportedSchemataShape

# -- "Global schema" -------------------------------------------------------
$globalSchemaShape = Get-ShapeByText $grp2.GroupItems "Global schema"
Remove-EndParaRPr $globalSchemaShape

# -- "(Domenig & Dittrich 1999 Sigmod Record)" -----------------------------
$refShape = $null
for ($i = 1; $i -le $slide2.Shapes.Count; $i++) {
    $sh = $slide2.Shapes.Item($i)
    if ($sh.HasTextFrame -and $sh.TextFrame.HasText) {
        if ($sh.TextFrame.TextRange.Text -like "*Domenig*") {
            $refShape = $sh
        }
    }
}
$refTr = $refShape.TextFrame2.TextRange
$fullText = $refTr.Text
$startIdx = $fullText.IndexOf("Domenig") + 1
$mergedLen = "Domenig & Dittrich 1999 Sigmod Record)".Length
$refMerged = $refTr.Characters($startIdx, $mergedLen)
$refMerged.Text = "Domenig & Dittrich 1999 Sigmod Record)"

# ----------------------------------------------------------------------
# Slide 3
# ----------------------------------------------------------------------
$slide3 = $p.Slides.Item(3)
$grp3 = $null
for ($i = 1; $i -le $slide3.Shapes.Count; $i++) {
    $sh = $slide3.Shapes.Item($i)
    if ($sh.Type -eq 6 -and $sh.Name -eq "Grupo 3") {
        $grp3 = $sh
    }
}

$distShape = Get-ShapeByText $grp3.GroupItems "Distributed data services"
Remove-EndParaRPr $distShape

$apiShape = Get-ShapeByText $grp3.GroupItems "Exported API"
Remove-EndParaRPr $apiShape

# -- drop the two click animations (Fade in on the two top-level groups) --
$seq = $slide3.TimeLine.MainSequence
while ($seq.Count -gt 0) {
    $seq.Item(1).Delete()
}
